$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark the rows that now have a completed Milestone II feature (roman numeral "II"
# in column E, and an "X" completion mark in column F). These drive the
# point-light / directional-light combo rows called out in the commit message.
$rows = @(5, 7, 9, 30, 40, 41, 46)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value = "II"
    $ws.Cells.Item($r, 6).Value = "X"
}

# Row 91 (Milestone II carry-over) is also marked complete with an X.
$ws.Range("D91").Value = "X"

# Leave the selection where the author last left it.
$ws.Range("F9").Select()
